$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A (Date) and D (Week) hold values that look like a date / plain
# number to the parser ("2023-05-28" and "21"). Pre-format as text and
# clear the formatting afterwards so the stored value keeps its literal
# string form while the cell ends up with the default (unstyled) look,
# matching the other data rows in the sheet.
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "2023-05-28"
$ws.Range("A4").ClearFormats()

$ws.Range("B4").Value = "15:06:41"
$ws.Range("C4").Value = "Sunday"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "21"
$ws.Range("D4").ClearFormats()

$ws.Range("E4").Value = 119926
$ws.Range("F4").Value = 133334
$ws.Range("G4").Value = 157695
$ws.Range("H4").Value = 130849
$ws.Range("I4").Value = 174265
$ws.Range("J4").Value = 114242
$ws.Range("K4").Value = 198116
$ws.Range("L4").Value = 219609
$ws.Range("M4").Value = 171953
$ws.Range("N4").Value = 119654
$ws.Range("O4").Value = 38541
$ws.Range("P4").Value = 34916
$ws.Range("Q4").Value = 50339
$ws.Range("R4").Value = -1
$ws.Range("S4").Value = 36871
$ws.Range("T4").Value = -1
